# Generate Report for Handoff
# Removes the "e79bfd45-f201-4ebb-a45c-88dc247ce92b" record (row 3) from every
# sheet, and refreshes the "7993d117-3123-48eb-930f-6b13ea96c3a6" record
# (row 2) to show it is now ready for handoff with an updated timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Remove stale hyperlinks that target row 3 on each sheet --------------
# (deleted one at a time via a fresh enumeration each time, since deleting
# mid-enumeration can skip entries as the collection re-indexes)

function Remove-HyperlinksAt($ws, $addresses) {
    foreach ($t in $addresses) {
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address() -eq $t) {
                $hl.Delete()
                break
            }
        }
    }
}

Remove-HyperlinksAt $overview @('$A$3')
Remove-HyperlinksAt $zhcn @('$A$3', '$B$3', '$D$3', '$F$3', '$G$3')
Remove-HyperlinksAt $dede @('$A$3', '$B$3', '$D$3', '$F$3', '$G$3')

# --- Delete row 3 (the e79bfd45... entry) on every sheet -------------------
$overview.Rows.Item(3).Delete()
$zhcn.Rows.Item(3).Delete()
$dede.Rows.Item(3).Delete()

# --- Refresh the remaining record (row 2) to reflect the handoff ----------
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-39-18 22:39:55"

$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-18 22:39:52"

$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-18 22:39:55"
